$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CDRDfRCP")

# Set clean firm requirement for CES ("required share peak") to 0
$ws.Range("B3").Value = 0

# Move active selection to B4 to match saved cursor position
$ws.Activate()
$ws.Range("B4").Select()
